# resize_table: shrink/restyle the two result tables on slide 1.
# (formerly "format_table" - renamed per commit "update function format_table to resize_table")

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# EMU per point, used to convert target column widths (expressed in EMU,
# as they appear in the OOXML <a:gridCol w="..."/>) into the point values
# the COM `Column.Width` property expects.
$EMU_PER_PT = 12700.0

function Resize-TableColumns($table, $widthsEmu) {
    for ($i = 1; $i -le $widthsEmu.Count; $i++) {
        $table.Columns.Item($i).Width = $widthsEmu[$i - 1] / $EMU_PER_PT
    }
}

# ---------------------------------------------------------------------
# Table 3 (first table, shape index 3): shrink columns, shrink all text
# to 10pt, and center the header row.
# ---------------------------------------------------------------------
$tbl1 = $s.Shapes.Item(3).Table

Resize-TableColumns $tbl1 @(598000, 1741000, 1487000, 598000, 471000, 1487000)

$rows1 = $tbl1.Rows.Count
$cols1 = $tbl1.Columns.Count

for ($r = 1; $r -le $rows1; $r++) {
    for ($c = 1; $c -le $cols1; $c++) {
        $cell = $tbl1.Cell($r, $c)
        $tr = $cell.Shape.TextFrame.TextRange
        $tr.Font.Size = 10
        if ($r -eq 1) {
            $tr.ParagraphFormat.Alignment = 2  # ppAlignCenter
        }
    }
}

# ---------------------------------------------------------------------
# Table 4 (second table, shape index 4): widen columns and center the
# header row only (font size in this table is left as-is).
# ---------------------------------------------------------------------
$tbl2 = $s.Shapes.Item(4).Table

Resize-TableColumns $tbl2 @(699600, 2071200, 1766400, 699600, 547200, 1766400)

$cols2 = $tbl2.Columns.Count
for ($c = 1; $c -le $cols2; $c++) {
    $cell = $tbl2.Cell(1, $c)
    $tr = $cell.Shape.TextFrame.TextRange
    $tr.ParagraphFormat.Alignment = 2  # ppAlignCenter
}
